# The sheet "guild" has a single column (A) of entries (shared strings).
# Two entries - "앰살" (row 35) and "되겐냐" (originally row 53) - were
# removed from the list, with the remaining entries shifting up to fill
# the gap (i.e. a normal row delete, not a value clear).
#
# Deleting row 35 first shifts "되겐냐" up from row 53 to row 52, so the
# second delete must target row 52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(35).Delete()
$ws.Rows(52).Delete()
